$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 362 (old rows 362:411 shift down to 365:414)
$ws.Rows("362:364").Insert()

# Populate the 3 newly inserted rows with fresh data (row 362, 363, 364)
$ws.Range("A362").Value2 = 10
$ws.Range("B362").Value2 = "Vega Modelo de Temuco"
$ws.Range("C362").Value2 = "La Araucanía"
$ws.Range("D362").Value2 = 44449
$ws.Range("E362").Value2 = 9
$ws.Range("F362").Value2 = 100114001
$ws.Range("G362").Value2 = "Papa"
$ws.Range("H362").Value2 = "Asterix"
$ws.Range("I362").Value2 = "1a (guarda)"
$ws.Range("J362").Value2 = 1100
$ws.Range("K362").Value2 = 8000
$ws.Range("L362").Value2 = 9000
$ws.Range("M362").Value2 = 8409
$ws.Range("N362").Value2 = "$/malla 25 kilos"
$ws.Range("O362").Value2 = "Provincia de Cautín"
$ws.Range("P362").Value2 = 336
$ws.Range("Q362").Value2 = 25
$ws.Range("R362").Value2 = "Hortaliza"

$ws.Range("A363").Value2 = 10
$ws.Range("B363").Value2 = "Vega Modelo de Temuco"
$ws.Range("C363").Value2 = "La Araucanía"
$ws.Range("D363").Value2 = 44449
$ws.Range("E363").Value2 = 9
$ws.Range("F363").Value2 = 100114001
$ws.Range("G363").Value2 = "Papa"
$ws.Range("H363").Value2 = "Asterix"
$ws.Range("I363").Value2 = "1a (guarda)"
$ws.Range("J363").Value2 = 1450
$ws.Range("K363").Value2 = 7000
$ws.Range("L363").Value2 = 7500
$ws.Range("M363").Value2 = 7224
$ws.Range("N363").Value2 = "$/saco 25 kilos"
$ws.Range("O363").Value2 = "Provincia de Cautín"
$ws.Range("P363").Value2 = 289
$ws.Range("Q363").Value2 = 25
$ws.Range("R363").Value2 = "Hortaliza"

$ws.Range("A364").Value2 = 10
$ws.Range("B364").Value2 = "Vega Modelo de Temuco"
$ws.Range("C364").Value2 = "La Araucanía"
$ws.Range("D364").Value2 = 44449
$ws.Range("E364").Value2 = 9
$ws.Range("F364").Value2 = 100114001
$ws.Range("G364").Value2 = "Papa"
$ws.Range("H364").Value2 = "Rosara"
$ws.Range("I364").Value2 = "1a (guarda)"
$ws.Range("J364").Value2 = 550
$ws.Range("K364").Value2 = 7000
$ws.Range("L364").Value2 = 7500
$ws.Range("M364").Value2 = 7227
$ws.Range("N364").Value2 = "$/saco 25 kilos"
$ws.Range("O364").Value2 = "Provincia de Cautín"
$ws.Range("P364").Value2 = 289
$ws.Range("Q364").Value2 = 25
$ws.Range("R364").Value2 = "Hortaliza"
